# Add new BOM rows to the "Misc" sheet describing the new faraday cage parts.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Misc")

$rows = @(
    @("M-Rivet-Nut", "Part", "Reference model for M4 rivet nut: https://uk.farnell.com/tr-fastenings/w-m4-rede-oe-st-kg/rivet-nut-steel-electro-zinc-cr3/dp/2983977", "WIP"),
    @("M-RP-CAM", "Part ", "Reference model for Raspberry Pi Camera Module 3: https://datasheets.raspberrypi.com/camera/camera-module-3-standard-mechanical-drawing.pdf", "WIP"),
    @("M-SmartDriveDuo-10", "Part", "Reference model for SmartDriveDuo-10: https://docs.google.com/document/d/1ECKsH1WfgZLomwFoPQ8PlkUzVEm88A6Sup9X0wDBHKA/view", "WIP"),
    @("M-TL-FC111B-20", "Part", "Reference model for TL-FC111B-20: https://www.tp-link.com/nl/business-networking/accessory/tl-fc111b-20/", "WIP"),
    @("M-12V-12Ah-Lead-Acid", "Part", "Reference model for Multipower 12V 12Ah lead acid (4.8mm): https://eu.nkon.nl/rechargeable/lead-acid-batteries/12v-vrla/multipower-12v-12ah-48.html", "WIP")
)

$startRow = 10
for ($i = 0; $i -lt $rows.Count; $i++) {
    $r = $startRow + $i
    $ws.Cells.Item($r, 1).Value = $rows[$i][0]
    $ws.Cells.Item($r, 2).Value = $rows[$i][1]
    $ws.Cells.Item($r, 3).Value = $rows[$i][2]
    $ws.Cells.Item($r, 4).Value = $rows[$i][3]
}

$ws.Columns.Item(1).AutoFit() | Out-Null
$ws.Columns.Item(3).AutoFit() | Out-Null
$ws.Range("C20").Select()
